$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split "gPOC/m/yr" / "gDOC/m/yr" into per-day rate units: "gPOC/m/d" / "gDOC/m/d"
$ws.Range("C14").Value = "gPOC/m/d"
$ws.Range("C15").Value = "gDOC/m/d"

# Update the window scroll/selection state: remove pinned top-left row, move selection to C16
$ws.Range("C16").Select()
